$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the workers/periods data table (rows 16-20, columns C=doc number, D=name, E=period)
# New data adds period 1711 for both workers and keeps all existing (worker, period)
# combinations, now grouped by worker and ordered by period.
$ws.Range("C16").Value = "73196852"
$ws.Range("D16").Value = "JORGE LUIS TAPIAS ROJAS"
$ws.Range("E16").Value = "1711"

$ws.Range("C17").Value = "9153002"
$ws.Range("D17").Value = "FELIX VALENCIA PEREZ"
$ws.Range("E17").Value = "1711"

$ws.Range("C18").Value = "73196852"
$ws.Range("D18").Value = "JORGE LUIS TAPIAS ROJAS"
$ws.Range("E18").Value = "1712"

$ws.Range("C19").Value = "9153002"
$ws.Range("D19").Value = "FELIX VALENCIA PEREZ"
$ws.Range("E19").Value = "1712"

$ws.Range("C20").Value = "9153002"
$ws.Range("D20").Value = "FELIX VALENCIA PEREZ"
$ws.Range("E20").Value = "1801"

$wb.Save()
